# Auto-generated edit script applying numeric updates to the 'Alpha_Profits' workbook
# (scheduled runner refresh of market-price-derived columns H-N across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2751.1667
$ws.Range("J17").Value = 2751.1667
$ws.Range("L17").Value = 8253.500100000001
$ws.Range("N17").Value = -8589.500100000001

$ws.Range("H28").Value = 1002.2857
$ws.Range("I28").Value = 1002.5
$ws.Range("K28").Value = 1002.5
$ws.Range("M28").Value = -517.5

$ws.Range("H33").Value = 34628.83
$ws.Range("I33").Value = 35863.465
$ws.Range("K33").Value = 35863.465
$ws.Range("M33").Value = -35634.465

$ws.Range("H40").Value = 3860.5557
$ws.Range("I40").Value = 3466.6667
$ws.Range("J40").Value = 4057.5
$ws.Range("K40").Value = 3466.6667
$ws.Range("L40").Value = 4057.5
$ws.Range("M40").Value = -3291.6667
$ws.Range("N40").Value = -4407.5

$ws.Range("H92").Value = 845.96155
$ws.Range("I92").Value = 965.8333
$ws.Range("K92").Value = 965.8333
$ws.Range("M92").Value = 282.1667

$ws.Range("H112").Value = 3728.3684
$ws.Range("J112").Value = 3746.8333
$ws.Range("L112").Value = 11240.4999
$ws.Range("N112").Value = -13456.4999

$ws.Range("H116").Value = 2990
$ws.Range("I116").Value = 2700
$ws.Range("J116").Value = 3425
$ws.Range("K116").Value = 2700
$ws.Range("L116").Value = 3425
$ws.Range("M116").Value = 742
$ws.Range("N116").Value = -10309

$ws.Range("H138").Value = 2672.12
$ws.Range("I138").Value = 1398.1
$ws.Range("J138").Value = 2990.625
$ws.Range("K138").Value = 4194.299999999999
$ws.Range("L138").Value = 8971.875
$ws.Range("M138").Value = 945.7000000000007
$ws.Range("N138").Value = -19251.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3005.6
$ws.Range("I61").Value = 2898.8462
$ws.Range("J61").Value = 3699.5
$ws.Range("K61").Value = 2898.8462
$ws.Range("L61").Value = 3699.5
$ws.Range("M61").Value = -2686.8462
$ws.Range("N61").Value = -4123.5

$ws.Range("H110").Value = 3359.4722
$ws.Range("I110").Value = 2247.3794
$ws.Range("J110").Value = 7966.7144
$ws.Range("K110").Value = 2247.3794
$ws.Range("L110").Value = 7966.7144
$ws.Range("M110").Value = -202.3793999999998
$ws.Range("N110").Value = -12056.7144

$ws.Range("H136").Value = 3005.6
$ws.Range("I136").Value = 2898.8462
$ws.Range("J136").Value = 3699.5
$ws.Range("K136").Value = 8696.5386
$ws.Range("L136").Value = 11098.5
$ws.Range("M136").Value = -6146.5386
$ws.Range("N136").Value = -16198.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3728.375
$ws.Range("J107").Value = 3603.6
$ws.Range("L107").Value = 3603.6
$ws.Range("N107").Value = -7443.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2781.2
$ws.Range("I16").Value = 2669.3333
$ws.Range("J16").Value = 2949
$ws.Range("K16").Value = 2669.3333
$ws.Range("L16").Value = 2949
$ws.Range("M16").Value = -2382.3333
$ws.Range("N16").Value = -3523

$ws.Range("H22").Value = 188.90909
$ws.Range("J22").Value = 220.5
$ws.Range("L22").Value = 220.5
$ws.Range("N22").Value = -920.5

$ws.Range("H70").Value = 35644.332
$ws.Range("J70").Value = 35644.332
$ws.Range("L70").Value = 35644.332
$ws.Range("N70").Value = -36274.332

$ws.Range("H73").Value = 35644.332
$ws.Range("J73").Value = 35644.332
$ws.Range("L73").Value = 35644.332
$ws.Range("N73").Value = -37828.332

$ws.Range("H113").Value = 2781.2
$ws.Range("I113").Value = 2669.3333
$ws.Range("J113").Value = 2949
$ws.Range("K113").Value = 2669.3333
$ws.Range("L113").Value = 2949
$ws.Range("M113").Value = -499.3332999999998
$ws.Range("N113").Value = -7289

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2991.6667
$ws.Range("J51").Value = 2990
$ws.Range("L51").Value = 8970
$ws.Range("N51").Value = -9890

$ws.Range("H57").Value = 6250
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 6250
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 18750
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -19868

$ws.Range("H58").Value = 1300
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 1450
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 4350
$ws.Range("N58").Value = -4606
$ws.Range("M58").Value = -2872

$ws.Range("H131").Value = 1791.7675
$ws.Range("I131").Value = 900
$ws.Range("J131").Value = 1835.2683
$ws.Range("K131").Value = 2700
$ws.Range("L131").Value = 5505.8049
$ws.Range("M131").Value = 2340
$ws.Range("N131").Value = -15585.8049

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2813.6428
$ws.Range("I113").Value = 2836
$ws.Range("K113").Value = 2836
$ws.Range("M113").Value = -666

$ws.Range("H132").Value = 1012.63336
$ws.Range("I132").Value = 882.2692
$ws.Range("J132").Value = 1860
$ws.Range("K132").Value = 2646.8076
$ws.Range("L132").Value = 5580
$ws.Range("M132").Value = -116.8076000000001
$ws.Range("N132").Value = -10640

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3991.889
$ws.Range("I7").Value = 3750.25
$ws.Range("K7").Value = 3750.25
$ws.Range("M7").Value = -3638.25

$ws.Range("H61").Value = 2609.5
$ws.Range("I61").Value = 2013.2858
$ws.Range("J61").Value = 4000.6667
$ws.Range("K61").Value = 2013.2858
$ws.Range("L61").Value = 4000.6667
$ws.Range("M61").Value = -1811.2858
$ws.Range("N61").Value = -4404.6667

$ws.Range("H93").Value = 15084.531
$ws.Range("I93").Value = 1757.421
$ws.Range("K93").Value = 1757.421
$ws.Range("M93").Value = -509.421

$ws.Range("H113").Value = 2609.5
$ws.Range("I113").Value = 2013.2858
$ws.Range("J113").Value = 4000.6667
$ws.Range("K113").Value = 2013.2858
$ws.Range("L113").Value = 4000.6667
$ws.Range("M113").Value = 156.7141999999999
$ws.Range("N113").Value = -8340.6667

$ws.Range("H126").Value = 3991.889
$ws.Range("I126").Value = 3750.25
$ws.Range("K126").Value = 11250.75
$ws.Range("M126").Value = -8780.75

$ws.Range("H132").Value = 2910.5833
$ws.Range("I132").Value = 2910.5833
$ws.Range("K132").Value = 8731.749899999999
$ws.Range("M132").Value = -6201.749899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1666
$ws.Range("I107").Value = 1531.1818
$ws.Range("K107").Value = 4593.5454
$ws.Range("M107").Value = -2673.5454

$ws.Range("H132").Value = 2028.8868
$ws.Range("I132").Value = 2176.9092
$ws.Range("K132").Value = 6530.7276
$ws.Range("M132").Value = -4000.7276
